$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column map: D=4 (Fecha), L=12 (Calidad), M=13 (Volumen), N=14 (Precio minimo),
# O=15 (Precio maximo), P=16 (Precio promedio ponderado), S=19 (Precio $/Kg)

$dateFmt = $ws.Cells.Item(2, 4).NumberFormat

# Row 43
$ws.Cells.Item(43, 4).Value = 44637
$ws.Cells.Item(43, 12).Value = "Especial"
$ws.Cells.Item(43, 13).Value = 240
$ws.Cells.Item(43, 14).Value = 12500
$ws.Cells.Item(43, 15).Value = 13000
$ws.Cells.Item(43, 16).Value = 12750
$ws.Cells.Item(43, 19).Value = 708

# Row 44
$ws.Cells.Item(44, 4).Value = 44637
$ws.Cells.Item(44, 12).Value = "Primera"
$ws.Cells.Item(44, 13).Value = 240
$ws.Cells.Item(44, 14).Value = 10500
$ws.Cells.Item(44, 15).Value = 11000
$ws.Cells.Item(44, 16).Value = 10750
$ws.Cells.Item(44, 19).Value = 597

# Row 45
$ws.Cells.Item(45, 4).Value = 44294
$ws.Cells.Item(45, 12).Value = "Especial"
$ws.Cells.Item(45, 13).Value = 360
$ws.Cells.Item(45, 14).Value = 12500
$ws.Cells.Item(45, 15).Value = 13000
$ws.Cells.Item(45, 16).Value = 12750
$ws.Cells.Item(45, 19).Value = 708

# Row 46
$ws.Cells.Item(46, 4).Value = 44294
$ws.Cells.Item(46, 12).Value = "Primera"
$ws.Cells.Item(46, 13).Value = 240
$ws.Cells.Item(46, 14).Value = 10500
$ws.Cells.Item(46, 15).Value = 11000
$ws.Cells.Item(46, 16).Value = 10750
$ws.Cells.Item(46, 19).Value = 597

# Row 47
$ws.Cells.Item(47, 4).Value = 44294
$ws.Cells.Item(47, 12).Value = "Segunda"
$ws.Cells.Item(47, 13).Value = 240
$ws.Cells.Item(47, 14).Value = 8500
$ws.Cells.Item(47, 15).Value = 9000
$ws.Cells.Item(47, 16).Value = 8750
$ws.Cells.Item(47, 19).Value = 486

# Row 48
$ws.Cells.Item(48, 4).Value = 44279
$ws.Cells.Item(48, 12).Value = "Especial"
$ws.Cells.Item(48, 13).Value = 200
$ws.Cells.Item(48, 14).Value = 12500
$ws.Cells.Item(48, 15).Value = 13000
$ws.Cells.Item(48, 16).Value = 12750
$ws.Cells.Item(48, 19).Value = 708

# Row 49
$ws.Cells.Item(49, 4).Value = 44279
$ws.Cells.Item(49, 12).Value = "Primera"
$ws.Cells.Item(49, 13).Value = 240
$ws.Cells.Item(49, 14).Value = 10500
$ws.Cells.Item(49, 15).Value = 11000
$ws.Cells.Item(49, 16).Value = 10750
$ws.Cells.Item(49, 19).Value = 597

# Row 50
$ws.Cells.Item(50, 4).Value = 44279
$ws.Cells.Item(50, 12).Value = "Segunda"
$ws.Cells.Item(50, 13).Value = 240
$ws.Cells.Item(50, 14).Value = 8500
$ws.Cells.Item(50, 15).Value = 9000
$ws.Cells.Item(50, 16).Value = 8750
$ws.Cells.Item(50, 19).Value = 486

# Row 51
$ws.Cells.Item(51, 4).Value = 44385
$ws.Cells.Item(51, 12).Value = "Especial"
$ws.Cells.Item(51, 13).Value = 120
$ws.Cells.Item(51, 14).Value = 14000
$ws.Cells.Item(51, 15).Value = 14500
$ws.Cells.Item(51, 16).Value = 14250
$ws.Cells.Item(51, 19).Value = 792

# Row 52
$ws.Cells.Item(52, 4).Value = 44385
$ws.Cells.Item(52, 12).Value = "Primera"
$ws.Cells.Item(52, 13).Value = 300
$ws.Cells.Item(52, 14).Value = 11000
$ws.Cells.Item(52, 15).Value = 11500
$ws.Cells.Item(52, 16).Value = 11250
$ws.Cells.Item(52, 19).Value = 625

# Row 53
$ws.Cells.Item(53, 4).Value = 44385
$ws.Cells.Item(53, 12).Value = "Segunda"
$ws.Cells.Item(53, 13).Value = 240
$ws.Cells.Item(53, 14).Value = 8000
$ws.Cells.Item(53, 15).Value = 8500
$ws.Cells.Item(53, 16).Value = 8250
$ws.Cells.Item(53, 19).Value = 458

# Row 54
$ws.Cells.Item(54, 4).Value = 44385
$ws.Cells.Item(54, 12).Value = "Tercera"
$ws.Cells.Item(54, 13).Value = 120
$ws.Cells.Item(54, 14).Value = 5000
$ws.Cells.Item(54, 15).Value = 5500
$ws.Cells.Item(54, 16).Value = 5250
$ws.Cells.Item(54, 19).Value = 292

# Row 55
$ws.Cells.Item(55, 4).Value = 44272
$ws.Cells.Item(55, 12).Value = "Especial"
$ws.Cells.Item(55, 13).Value = 160
$ws.Cells.Item(55, 14).Value = 12500
$ws.Cells.Item(55, 15).Value = 13000
$ws.Cells.Item(55, 16).Value = 12750
$ws.Cells.Item(55, 19).Value = 708

# Row 56
$ws.Cells.Item(56, 4).Value = 44272
$ws.Cells.Item(56, 12).Value = "Primera"
$ws.Cells.Item(56, 13).Value = 300
$ws.Cells.Item(56, 14).Value = 10500
$ws.Cells.Item(56, 15).Value = 11000
$ws.Cells.Item(56, 16).Value = 10750
$ws.Cells.Item(56, 19).Value = 597

# Row 57
$ws.Cells.Item(57, 4).Value = 44272
$ws.Cells.Item(57, 12).Value = "Segunda"
$ws.Cells.Item(57, 13).Value = 240
$ws.Cells.Item(57, 14).Value = 8500
$ws.Cells.Item(57, 15).Value = 9000
$ws.Cells.Item(57, 16).Value = 8750
$ws.Cells.Item(57, 19).Value = 486

# Row 58
$ws.Cells.Item(58, 4).Value = 44615
$ws.Cells.Item(58, 12).Value = "Especial"
$ws.Cells.Item(58, 13).Value = 200
$ws.Cells.Item(58, 14).Value = 14000
$ws.Cells.Item(58, 15).Value = 15000
$ws.Cells.Item(58, 16).Value = 14500
$ws.Cells.Item(58, 19).Value = 806

# Row 59
$ws.Cells.Item(59, 4).Value = 44615
$ws.Cells.Item(59, 12).Value = "Primera"
$ws.Cells.Item(59, 13).Value = 400
$ws.Cells.Item(59, 14).Value = 12000
$ws.Cells.Item(59, 15).Value = 13000
$ws.Cells.Item(59, 16).Value = 12500
$ws.Cells.Item(59, 19).Value = 694

# Row 60
$ws.Cells.Item(60, 4).Value = 44335
$ws.Cells.Item(60, 12).Value = "Especial"
$ws.Cells.Item(60, 13).Value = 240
$ws.Cells.Item(60, 14).Value = 19500
$ws.Cells.Item(60, 15).Value = 20000
$ws.Cells.Item(60, 16).Value = 19750
$ws.Cells.Item(60, 19).Value = 1097

# Row 61
$ws.Cells.Item(61, 1).Value = 2
$ws.Cells.Item(61, 2).Value = "Comercializadora del Agro de Limarí"
$ws.Cells.Item(61, 3).Value = "Coquimbo"
$ws.Cells.Item(61, 5).Value = 4
$ws.Cells.Item(61, 6).Value = "Fruta"
$ws.Cells.Item(61, 7).Value = 100107
$ws.Cells.Item(61, 8).Value = "Otros"
$ws.Cells.Item(61, 9).Value = 100107011
$ws.Cells.Item(61, 10).Value = "Tuna"
$ws.Cells.Item(61, 11).Value = "Sin especificar"
$ws.Cells.Item(61, 17).Value = "$/caja 18 kilos"
$ws.Cells.Item(61, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(61, 20).Value = 18
$ws.Cells.Item(61, 4).Value = 44335
$ws.Cells.Item(61, 4).NumberFormat = $dateFmt
$ws.Cells.Item(61, 12).Value = "Primera"
$ws.Cells.Item(61, 13).Value = 200
$ws.Cells.Item(61, 14).Value = 17500
$ws.Cells.Item(61, 15).Value = 18000
$ws.Cells.Item(61, 16).Value = 17750
$ws.Cells.Item(61, 19).Value = 986

# Row 62
$ws.Cells.Item(62, 1).Value = 2
$ws.Cells.Item(62, 2).Value = "Comercializadora del Agro de Limarí"
$ws.Cells.Item(62, 3).Value = "Coquimbo"
$ws.Cells.Item(62, 5).Value = 4
$ws.Cells.Item(62, 6).Value = "Fruta"
$ws.Cells.Item(62, 7).Value = 100107
$ws.Cells.Item(62, 8).Value = "Otros"
$ws.Cells.Item(62, 9).Value = 100107011
$ws.Cells.Item(62, 10).Value = "Tuna"
$ws.Cells.Item(62, 11).Value = "Sin especificar"
$ws.Cells.Item(62, 17).Value = "$/caja 18 kilos"
$ws.Cells.Item(62, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(62, 20).Value = 18
$ws.Cells.Item(62, 4).Value = 44335
$ws.Cells.Item(62, 4).NumberFormat = $dateFmt
$ws.Cells.Item(62, 12).Value = "Segunda"
$ws.Cells.Item(62, 13).Value = 160
$ws.Cells.Item(62, 14).Value = 12500
$ws.Cells.Item(62, 15).Value = 13000
$ws.Cells.Item(62, 16).Value = 12750
$ws.Cells.Item(62, 19).Value = 708
